# Yearly-Report-2021-DE456232.xlsx - Automation HUB update
# Rewrites the changed cells of the data table on Sheet1 (rows 3-8) to
# reflect the refreshed export. Values that look numeric or date-like are
# written with a temporary "@" (text) number format so Excel keeps storing
# them as text (shared strings) instead of auto-converting them to
# numbers/dates, matching the source data which is plain text. The number
# format is cleared again right after the write so the cell keeps its
# original (default) style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = ""
}

# Row 3: was 197993/Professional Services/2017-09-12/65658/13131.6/78789.6/USD
Set-TextValue "A3" "857700"
Set-TextValue "C3" "2017-07-01"
Set-TextValue "D3" "294833"
Set-TextValue "E3" "58966.6"
Set-TextValue "F3" "353800"
$ws.Range("G3").Value = "RON"

# Row 4: was 936027/Concierge Services/2017-06-12/288581/57716.2/346297/RON
Set-TextValue "A4" "857700"
$ws.Range("B4").Value = "Professional Services"
Set-TextValue "C4" "2017-07-01"
Set-TextValue "D4" "294833"
Set-TextValue "E4" "58966.6"
Set-TextValue "F4" "353800"

# Row 5: was 571895/Concierge Services/2017-06-07/293171/58634.2/351805/EUR
Set-TextValue "A5" "197993"
$ws.Range("B5").Value = "Professional Services"
Set-TextValue "C5" "2017-09-12"
Set-TextValue "D5" "65658"
Set-TextValue "E5" "13131.6"
Set-TextValue "F5" "78789.6"
$ws.Range("G5").Value = "USD"

# Row 6: was 197993/Professional Services/2017-09-12/65658/13131.6/78789.6/USD
Set-TextValue "A6" "185926"
$ws.Range("B6").Value = "Various paper supplies"
Set-TextValue "C6" "2017-10-21"
Set-TextValue "D6" "107274"
Set-TextValue "E6" "21454.8"
Set-TextValue "F6" "128729"

# Row 7: was 197993/Professional Services/2017-09-12/65658/13131.6/78789.6/USD
Set-TextValue "A7" "783399"
$ws.Range("B7").Value = "Beverages and Catering"
Set-TextValue "C7" "2017-10-16"
Set-TextValue "D7" "21436"
Set-TextValue "E7" "4287.2"
Set-TextValue "F7" "25723.2"

# Row 8: was 197993/Professional Services/2017-09-12/65658/13131.6/78789.6/USD
Set-TextValue "A8" "783399"
$ws.Range("B8").Value = "Beverages and Catering"
Set-TextValue "C8" "2017-10-16"
Set-TextValue "D8" "21436"
Set-TextValue "E8" "4287.2"
Set-TextValue "F8" "25723.2"
